# Generate Report for Handoff
# Rename the source GUID-based file references from
#   ed68e2ac-d253-4a56-b399-99ce449a5049
# to
#   bc7ad01a-d8c4-4593-bace-17fb2811f112
# and bump the related timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "ed68e2ac-d253-4a56-b399-99ce449a5049"
$newGuid = "bc7ad01a-d8c4-4593-bace-17fb2811f112"

$oldZhHash = "5cc561ff277242e81f731fea8121977fad73065a"
$newZhHash = "531a9394bd1e7a4793c0429ba8ba9aa0cc169170"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-08-13 19:12:45"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Range("G2").Value = "$newGuid.$newZhHash.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-13 19:12:37"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Range("G2").Value = "$newGuid.$newZhHash.de-de.xlf"
